$wb = $excel.ActiveWorkbook

# --- Update "Score" sheet (sheet1) with recalculated scenario values ---
$wsScore = $wb.Worksheets.Item("Score")

$wsScore.Range("D4").Value2 = 0.4431831111395371
$wsScore.Range("E4").Value2 = 0.6209707718945657
$wsScore.Range("F4").Value2 = 0.4216145211398413
$wsScore.Range("D5").Value2 = 0.7523279998200449
$wsScore.Range("D6").Value2 = 0.383492181287474
$wsScore.Range("D7").Value2 = 0.6499145417213249
$wsScore.Range("E7").Value2 = 0.8778272545851326
$wsScore.Range("F7").Value2 = 0.1323819684936061
$wsScore.Range("D9").Value2 = 0.3620185166548159
$wsScore.Range("E9").Value2 = 0.6154564462793018
$wsScore.Range("F9").Value2 = 0.4122532093345651
$wsScore.Range("D10").Value2 = 0.6921888206727864
$wsScore.Range("E10").Value2 = 0.8737372697874098
$wsScore.Range("F10").Value2 = 0.136786239334902
$wsScore.Range("D11").Value2 = 0.4721356533617872
$wsScore.Range("E11").Value2 = 0.5857689352793057
$wsScore.Range("F11").Value2 = 0.441314105938786
$wsScore.Range("D12").Value2 = 0.7447688187902374
$wsScore.Range("E12").Value2 = 0.9259958224920952
$wsScore.Range("F12").Value2 = 0.07400417750790486
$wsScore.Range("D13").Value2 = 0.4046322827814371
$wsScore.Range("E13").Value2 = 0.789637417678884
$wsScore.Range("F13").Value2 = 0.2103625823211161
$wsScore.Range("D14").Value2 = 0.3620185166548159
$wsScore.Range("E14").Value2 = 0.6154564462793018
$wsScore.Range("F14").Value2 = 0.4122532093345651
$wsScore.Range("D15").Value2 = 0.6921888206727864
$wsScore.Range("E15").Value2 = 0.8737372697874098
$wsScore.Range("F15").Value2 = 0.136786239334902
$wsScore.Range("D16").Value2 = 0.6729456182270112
$wsScore.Range("E16").Value2 = 0.8954214936130447
$wsScore.Range("F16").Value2 = 0.1147919572510599
$wsScore.Range("D17").Value2 = 0.4126776354512215
$wsScore.Range("E17").Value2 = 0.7329567446632728
$wsScore.Range("F17").Value2 = 0.2883098614051222
$wsScore.Range("D18").Value2 = 0.652467036162367
$wsScore.Range("E18").Value2 = 0.8948223230462478
$wsScore.Range("F18").Value2 = 0.1175278946680752
$wsScore.Range("D19").Value2 = 0.6763014222131597
$wsScore.Range("E19").Value2 = 0.9182816825663566
$wsScore.Range("F19").Value2 = 0.09536521761552316
$wsScore.Range("D20").Value2 = 0.6872596566604285
$wsScore.Range("E20").Value2 = 0.9945488454708578
$wsScore.Range("F20").Value2 = 0.005451154529142241
$wsScore.Range("D21").Value2 = 0.357159580321769
$wsScore.Range("E21").Value2 = 0.7450870790206342
$wsScore.Range("F21").Value2 = 0.2811931101015958
$wsScore.Range("D22").Value2 = 0.3141296112295865
$wsScore.Range("E22").Value2 = 0.6405381660692044
$wsScore.Range("F22").Value2 = 0.3952126938577485
$wsScore.Range("D23").Value2 = 0.3977658100762115
$wsScore.Range("E23").Value2 = 0.7779344526405685
$wsScore.Range("F23").Value2 = 0.2529864858380443
$wsScore.Range("D24").Value2 = 0.7328119078577248
$wsScore.Range("E24").Value2 = 0.9154102190182059
$wsScore.Range("F24").Value2 = 0.09871345336491412
$wsScore.Range("D25").Value2 = 0.6611790462448828
$wsScore.Range("E25").Value2 = 0.9268056466598925
$wsScore.Range("F25").Value2 = 0.07735996034522688
$wsScore.Range("D26").Value2 = 0.7103339116391424
$wsScore.Range("E26").Value2 = 0.9242508536484528
$wsScore.Range("F26").Value2 = 0.08007335687116421
$wsScore.Range("D27").Value2 = 0.6965069871685858
$wsScore.Range("E27").Value2 = 0.8912367314413862
$wsScore.Range("F27").Value2 = 0.1215120652142653
$wsScore.Range("D29").Value2 = 0.4992126290541916
$wsScore.Range("E29").Value2 = 0.5487710217178979
$wsScore.Range("F29").Value2 = 0.4642226902831631
$wsScore.Range("D30").Value2 = 0.7254007048602344
$wsScore.Range("E30").Value2 = 0.8633910668752239
$wsScore.Range("F30").Value2 = 0.136608933124776
$wsScore.Range("D31").Value2 = 0.441158559173889
$wsScore.Range("E31").Value2 = 0.6755971720939401
$wsScore.Range("F31").Value2 = 0.32440282790606
$wsScore.Range("E32").Value2 = 0.6443363925673763
$wsScore.Range("D33").Value2 = 0.2184711574160119
$wsScore.Range("E33").Value2 = 0.6526356209817371
$wsScore.Range("F33").Value2 = 0.3769164705897979
$wsScore.Range("D34").Value2 = 0.6318552249037735
$wsScore.Range("E34").Value2 = 0.8635723314934812
$wsScore.Range("F34").Value2 = 0.144065893256218
$wsScore.Range("D35").Value2 = 0.266394388412658
$wsScore.Range("E35").Value2 = 0.7234538650328752
$wsScore.Range("F35").Value2 = 0.2938167438742204
$wsScore.Range("D36").Value2 = 0.4530893101245041
$wsScore.Range("E36").Value2 = 0.6093573956230762
$wsScore.Range("F36").Value2 = 0.4278797098413298
$wsScore.Range("D37").Value2 = 0.3920832712744233
$wsScore.Range("E37").Value2 = 0.628720671515097
$wsScore.Range("F37").Value2 = 0.4115790210535193
$wsScore.Range("D38").Value2 = 0.6899563758316446
$wsScore.Range("E38").Value2 = 0.8483425569755983
$wsScore.Range("F38").Value2 = 0.1749101136239687
$wsScore.Range("D39").Value2 = 0.4166663804362242
$wsScore.Range("E39").Value2 = 0.6928387995506464
$wsScore.Range("F39").Value2 = 0.3453114187595375
$wsScore.Range("D40").Value2 = 0.7503764233194545
$wsScore.Range("E40").Value2 = 0.9614811330946759
$wsScore.Range("F40").Value2 = 0.0385188669053242
$wsScore.Range("D41").Value2 = 0.7404550478014275
$wsScore.Range("E41").Value2 = 0.9602675573651627
$wsScore.Range("F41").Value2 = 0.04205949565625011
$wsScore.Range("D42").Value2 = 0.7470939785260781
$wsScore.Range("E42").Value2 = 0.9554698424335188
$wsScore.Range("F42").Value2 = 0.05234118452371551
$wsScore.Range("D43").Value2 = 0.7321698097688492
$wsScore.Range("E43").Value2 = 0.9971839665237002
$wsScore.Range("F43").Value2 = 0.002816033476299789
$wsScore.Range("D44").Value2 = 0.3894949150663777
$wsScore.Range("D45").Value2 = 0.3691789964565897
$wsScore.Range("E45").Value2 = 0.8778824806841222
$wsScore.Range("D46").Value2 = 0.3878507277490345
$wsScore.Range("E46").Value2 = 0.8643685873823383
$wsScore.Range("F46").Value2 = 0.1565802380758488
$wsScore.Range("D47").Value2 = 0.5956039862884859
$wsScore.Range("E47").Value2 = 0.9902242380567424
$wsScore.Range("F47").Value2 = 0.009775761943257579
$wsScore.Range("D48").Value2 = 0.5699763760797616
$wsScore.Range("E48").Value2 = 0.8548791716110304
$wsScore.Range("F48").Value2 = 0.1652951900634266
$wsScore.Range("D49").Value2 = 0.6926403650370847
$wsScore.Range("E49").Value2 = 0.9150088116147705
$wsScore.Range("F49").Value2 = 0.0921833286296143
$wsScore.Range("D50").Value2 = 0.5860524714493888
$wsScore.Range("E50").Value2 = 0.8191065114949702
$wsScore.Range("F50").Value2 = 0.200665423617032
$wsScore.Range("D51").Value2 = 0.7069724604640318
$wsScore.Range("E51").Value2 = 0.9092479535890007
$wsScore.Range("F51").Value2 = 0.1042389039145356

# --- Update "Rank" sheet (sheet2) with recalculated ranks ---
$wsRank = $wb.Worksheets.Item("Rank")

$wsRank.Range("D9").Value2 = 3
$wsRank.Range("F9").Value2 = 2
$wsRank.Range("D14").Value2 = 3
$wsRank.Range("F14").Value2 = 2
$wsRank.Range("D22").Value2 = 3
$wsRank.Range("F22").Value2 = 2
$wsRank.Range("D37").Value2 = 3
$wsRank.Range("F37").Value2 = 2
